$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.63908767700195
$ws.Range("C3").Value = 17.25101470947266
$ws.Range("C4").Value = 17.14611053466797
$ws.Range("C5").Value = 16.89291000366211
$ws.Range("C6").Value = 17.42100715637207
